$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '29.152.80'
Set-TextValue $ws.Range("E2") '  +0.03%  '

Set-TextValue $ws.Range("D3") '1.827.79'
Set-TextValue $ws.Range("E3") '  -0.55%  '

Set-TextValue $ws.Range("D4") '0.9990'
Set-TextValue $ws.Range("E4") '  -0.10%  '

Set-TextValue $ws.Range("D5") '242.14'
Set-TextValue $ws.Range("E5") '  -0.59%  '

Set-TextValue $ws.Range("E6") '  -0.94%  '

Set-TextValue $ws.Range("E7") '  -0.09%  '

Set-TextValue $ws.Range("D8") '0.07361'
Set-TextValue $ws.Range("E8") '  -2.07%  '

Set-TextValue $ws.Range("D9") '0.2913'
Set-TextValue $ws.Range("E9") '  -0.97%  '

Set-TextValue $ws.Range("D10") '23.12'
Set-TextValue $ws.Range("E10") '  -0.80%  '

Set-TextValue $ws.Range("D11") '0.07666'
Set-TextValue $ws.Range("E11") '  -0.51%  '

Set-TextValue $ws.Range("D12") '1.830.68'
Set-TextValue $ws.Range("E12") '  -0.55%  '

Set-TextValue $ws.Range("D13") '4.946'

Set-TextValue $ws.Range("D14") '0.6645'
Set-TextValue $ws.Range("E14") '  -1.78%  '

Set-TextValue $ws.Range("D15") '82.15'
Set-TextValue $ws.Range("E15") '  -1.20%  '

Set-TextValue $ws.Range("D16") '0.000008953'
Set-TextValue $ws.Range("E16") '  -3.54%  '

Set-TextValue $ws.Range("D17") '5.833'
Set-TextValue $ws.Range("E17") '  -2.23%  '

Set-TextValue $ws.Range("D18") '29.126.41'

Set-TextValue $ws.Range("D19") '2.074.59'
Set-TextValue $ws.Range("E19") '  -0.60%  '

Set-TextValue $ws.Range("D20") '236.38'
Set-TextValue $ws.Range("E20") '  +1.53%  '

Set-TextValue $ws.Range("D21") '12.42'
Set-TextValue $ws.Range("E21") '  -2.11%  '

Set-TextValue $ws.Range("D22") '0.9996'
Set-TextValue $ws.Range("E22") '  -0.16%  '

Set-TextValue $ws.Range("D23") '7.315'

Set-TextValue $ws.Range("D24") '1.000'
Set-TextValue $ws.Range("E24") '  -0.10%  '

Set-TextValue $ws.Range("E25") '  -1.32%  '

Set-TextValue $ws.Range("D26") '0.1406'
Set-TextValue $ws.Range("E26") '  +0.19%  '

Set-TextValue $ws.Range("D27") '8.495'
Set-TextValue $ws.Range("E27") '  -0.60%  '

Set-TextValue $ws.Range("E28") '  -1.79%  '

Set-TextValue $ws.Range("D29") '1.483'

Set-TextValue $ws.Range("D30") '0.05937'
Set-TextValue $ws.Range("E30") '  +6.44%  '

Set-TextValue $ws.Range("D31") '4.068'
Set-TextValue $ws.Range("E31") '  -2.05%  '

Set-TextValue $ws.Range("D32") '4.084'
Set-TextValue $ws.Range("E32") '  -2.48%  '

Set-TextValue $ws.Range("D33") '1.204'
Set-TextValue $ws.Range("E33") '  -0.45%  '

Set-TextValue $ws.Range("D34") '1.857'
Set-TextValue $ws.Range("E34") '  +0.07%  '

Set-TextValue $ws.Range("B35") 'ARBITRUM'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D35") '1.138'
Set-TextValue $ws.Range("E35") '  -0.85%  '

Set-TextValue $ws.Range("B36") 'ImmutableX'
Set-TextValue $ws.Range("C36") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D36") '0.7223'
Set-TextValue $ws.Range("E36") '  -3.76%  '

Set-TextValue $ws.Range("D37") '2.607'
Set-TextValue $ws.Range("E37") '  -1.97%  '

Set-TextValue $ws.Range("D38") '2.837'
Set-TextValue $ws.Range("E38") '  +2.41%  '

Set-TextValue $ws.Range("D39") '1.219.22'
Set-TextValue $ws.Range("E39") '  -1.95%  '

Set-TextValue $ws.Range("D40") '0.01751'
Set-TextValue $ws.Range("E40") '  -2.04%  '

Set-TextValue $ws.Range("D41") '0.9173'
Set-TextValue $ws.Range("E41") '  +1.92%  '

Set-TextValue $ws.Range("D42") '6.266'
Set-TextValue $ws.Range("E42") '  -5.15%  '

Set-TextValue $ws.Range("D43") '0.9997'
Set-TextValue $ws.Range("E43") '  -0.05%  '

Set-TextValue $ws.Range("D44") '101.76'
Set-TextValue $ws.Range("E44") '  -0.57%  '

Set-TextValue $ws.Range("D45") '1.980.30'
Set-TextValue $ws.Range("E45") '  -0.72%  '

Set-TextValue $ws.Range("D46") '64.73'
Set-TextValue $ws.Range("E46") '  -2.94%  '

Set-TextValue $ws.Range("D47") '0.5052'
Set-TextValue $ws.Range("E47") '  -0.74%  '

Set-TextValue $ws.Range("E48") '  -2.47%  '

Set-TextValue $ws.Range("D49") '0.4015'
Set-TextValue $ws.Range("E49") '  -1.83%  '

Set-TextValue $ws.Range("D50") '9.066'
Set-TextValue $ws.Range("E50") '  -0.24%  '

Set-TextValue $ws.Range("D51") '0.1132'
Set-TextValue $ws.Range("E51") '  +2.05%  '
